$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value2 = 24562.223
$ws.Cells.Item(17, 10).Value2 = 24562.223
$ws.Cells.Item(17, 12).Value2 = 73686.66900000001
$ws.Cells.Item(17, 14).Value2 = -74022.66900000001

# ALC row 108
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(108, 8).Value2 = 29553
$ws.Cells.Item(108, 10).Value2 = 29553
$ws.Cells.Item(108, 12).Value2 = 29553
$ws.Cells.Item(108, 14).Value2 = -37233

# ALC row 109
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(109, 8).Value2 = 37453.332
$ws.Cells.Item(109, 10).Value2 = 37453.332
$ws.Cells.Item(109, 12).Value2 = 37453.332
$ws.Cells.Item(109, 14).Value2 = -40227.332

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value2 = 9434.388999999999
$ws.Cells.Item(116, 9).Value2 = 8000
$ws.Cells.Item(116, 10).Value2 = 9721.267
$ws.Cells.Item(116, 11).Value2 = 8000
$ws.Cells.Item(116, 12).Value2 = 9721.267
$ws.Cells.Item(116, 13).Value2 = -4558
$ws.Cells.Item(116, 14).Value2 = -16605.267

# ALC row 117
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(117, 8).Value2 = 48542
$ws.Cells.Item(117, 10).Value2 = 48542
$ws.Cells.Item(117, 12).Value2 = 48542
$ws.Cells.Item(117, 14).Value2 = -57720

# ALC row 120
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(120, 8).Value2 = 48251.25
$ws.Cells.Item(120, 10).Value2 = 48251.25
$ws.Cells.Item(120, 12).Value2 = 48251.25
$ws.Cells.Item(120, 14).Value2 = -57927.25

# ALC row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(126, 8).Value2 = 39777.5
$ws.Cells.Item(126, 10).Value2 = 39777.5
$ws.Cells.Item(126, 12).Value2 = 39777.5
$ws.Cells.Item(126, 14).Value2 = -49657.5

# ALC row 128
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(128, 8).Value2 = 59979.332
$ws.Cells.Item(128, 10).Value2 = 59979.332
$ws.Cells.Item(128, 12).Value2 = 59979.332
$ws.Cells.Item(128, 14).Value2 = -69939.33199999999

# ALC row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(130, 8).Value2 = 45753.332
$ws.Cells.Item(130, 10).Value2 = 45753.332
$ws.Cells.Item(130, 12).Value2 = 45753.332
$ws.Cells.Item(130, 14).Value2 = -55793.332

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value2 = 2415.75
$ws.Cells.Item(138, 9).Value2 = 2362.2307
$ws.Cells.Item(138, 10).Value2 = 2439.7415
$ws.Cells.Item(138, 11).Value2 = 7086.6921
$ws.Cells.Item(138, 12).Value2 = 7319.2245
$ws.Cells.Item(138, 13).Value2 = -1946.6921
$ws.Cells.Item(138, 14).Value2 = -17599.2245

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value2 = 5103.3335
$ws.Cells.Item(141, 9).Value2 = 3123.5
$ws.Cells.Item(141, 10).Value2 = 9063
$ws.Cells.Item(141, 11).Value2 = 9370.5
$ws.Cells.Item(141, 12).Value2 = 27189
$ws.Cells.Item(141, 13).Value2 = -4190.5
$ws.Cells.Item(141, 14).Value2 = -37549

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value2 = 2264.3125
$ws.Cells.Item(45, 9).Value2 = 1936.4546
$ws.Cells.Item(45, 10).Value2 = 2985.6
$ws.Cells.Item(45, 11).Value2 = 1936.4546
$ws.Cells.Item(45, 12).Value2 = 2985.6
$ws.Cells.Item(45, 13).Value2 = -1559.4546
$ws.Cells.Item(45, 14).Value2 = -3739.6

# ARM row 107
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(107, 8).Value2 = 38500
$ws.Cells.Item(107, 10).Value2 = 38500
$ws.Cells.Item(107, 12).Value2 = 38500
$ws.Cells.Item(107, 14).Value2 = -46180

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(109, 8).Value2 = 41877
$ws.Cells.Item(109, 10).Value2 = 41877
$ws.Cells.Item(109, 12).Value2 = 41877
$ws.Cells.Item(109, 14).Value2 = -44651

# ARM row 117
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(117, 8).Value2 = 45698.8
$ws.Cells.Item(117, 10).Value2 = 45698.8
$ws.Cells.Item(117, 12).Value2 = 45698.8
$ws.Cells.Item(117, 14).Value2 = -54876.8

# ARM row 118
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(118, 8).Value2 = 46666.668
$ws.Cells.Item(118, 10).Value2 = 46666.668
$ws.Cells.Item(118, 12).Value2 = 46666.668
$ws.Cells.Item(118, 14).Value2 = -49980.668

# ARM row 120
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(120, 8).Value2 = 45346.668
$ws.Cells.Item(120, 10).Value2 = 45346.668
$ws.Cells.Item(120, 12).Value2 = 45346.668
$ws.Cells.Item(120, 14).Value2 = -55022.668

# ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value2 = 51429
$ws.Cells.Item(123, 9).Value2 = 0
$ws.Cells.Item(123, 11).Value2 = 0
$ws.Cells.Item(123, 13).ClearContents()

# ARM row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(125, 8).Value2 = 50707
$ws.Cells.Item(125, 10).Value2 = 50707
$ws.Cells.Item(125, 12).Value2 = 50707
$ws.Cells.Item(125, 14).Value2 = -60547

# ARM row 128
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(128, 8).Value2 = 50421
$ws.Cells.Item(128, 10).Value2 = 50421
$ws.Cells.Item(128, 12).Value2 = 50421
$ws.Cells.Item(128, 14).Value2 = -60381

# ARM row 131
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(131, 8).Value2 = 48681.668
$ws.Cells.Item(131, 10).Value2 = 48681.668
$ws.Cells.Item(131, 12).Value2 = 48681.668
$ws.Cells.Item(131, 14).Value2 = -58761.668

# BSM row 119
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(119, 8).Value2 = 48261
$ws.Cells.Item(119, 10).Value2 = 48261
$ws.Cells.Item(119, 12).Value2 = 48261
$ws.Cells.Item(119, 14).Value2 = -57937

# BSM row 120
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(120, 8).Value2 = 43394.332
$ws.Cells.Item(120, 10).Value2 = 43394.332
$ws.Cells.Item(120, 12).Value2 = 43394.332
$ws.Cells.Item(120, 14).Value2 = -53070.332

# BSM row 125
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(125, 8).Value2 = 50730
$ws.Cells.Item(125, 10).Value2 = 50730
$ws.Cells.Item(125, 12).Value2 = 50730
$ws.Cells.Item(125, 14).Value2 = -60570

# BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(126, 8).Value2 = 50772
$ws.Cells.Item(126, 10).Value2 = 50772
$ws.Cells.Item(126, 12).Value2 = 50772
$ws.Cells.Item(126, 14).Value2 = -60652

# CRP row 20
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value2 = 40283.168
$ws.Cells.Item(20, 10).Value2 = 40283.168
$ws.Cells.Item(20, 12).Value2 = 40283.168
$ws.Cells.Item(20, 14).Value2 = -40755.168

# CRP row 30
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(30, 8).Value2 = 40283.168
$ws.Cells.Item(30, 10).Value2 = 40283.168
$ws.Cells.Item(30, 12).Value2 = 40283.168
$ws.Cells.Item(30, 14).Value2 = -40465.168

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 4874.4263
$ws.Cells.Item(31, 9).Value2 = 1908.4814
$ws.Cells.Item(31, 10).Value2 = 6827.61
$ws.Cells.Item(31, 11).Value2 = 1908.4814
$ws.Cells.Item(31, 12).Value2 = 6827.61
$ws.Cells.Item(31, 13).Value2 = -1613.4814
$ws.Cells.Item(31, 14).Value2 = -7417.61

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value2 = 4874.4263
$ws.Cells.Item(34, 9).Value2 = 1908.4814
$ws.Cells.Item(34, 10).Value2 = 6827.61
$ws.Cells.Item(34, 11).Value2 = 1908.4814
$ws.Cells.Item(34, 12).Value2 = 6827.61
$ws.Cells.Item(34, 13).Value2 = -1706.4814
$ws.Cells.Item(34, 14).Value2 = -7231.61

# CRP row 100
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(100, 8).Value2 = 46996
$ws.Cells.Item(100, 10).Value2 = 46996
$ws.Cells.Item(100, 12).Value2 = 46996
$ws.Cells.Item(100, 14).Value2 = -49160

# CRP row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(116, 8).Value2 = 44078
$ws.Cells.Item(116, 10).Value2 = 44078
$ws.Cells.Item(116, 12).Value2 = 44078
$ws.Cells.Item(116, 14).Value2 = -53256

# CRP row 128
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(128, 8).Value2 = 40283.168
$ws.Cells.Item(128, 10).Value2 = 40283.168
$ws.Cells.Item(128, 12).Value2 = 40283.168
$ws.Cells.Item(128, 14).Value2 = -50243.168

# CUL row 50
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value2 = 166667200
$ws.Cells.Item(50, 9).Value2 = 76.666664
$ws.Cells.Item(50, 10).Value2 = 333334340
$ws.Cells.Item(50, 11).Value2 = 229.999992
$ws.Cells.Item(50, 12).Value2 = 1000003020
$ws.Cells.Item(50, 13).Value2 = 251.000008
$ws.Cells.Item(50, 14).Value2 = -1000003982

# CUL row 53
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(53, 8).Value2 = 166667200
$ws.Cells.Item(53, 9).Value2 = 76.666664
$ws.Cells.Item(53, 10).Value2 = 333334340
$ws.Cells.Item(53, 11).Value2 = 229.999992
$ws.Cells.Item(53, 12).Value2 = 1000003020
$ws.Cells.Item(53, 13).Value2 = 251.000008
$ws.Cells.Item(53, 14).Value2 = -1000003982

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value2 = 533.3333
$ws.Cells.Item(68, 9).Value2 = 320
$ws.Cells.Item(68, 10).Value2 = 800
$ws.Cells.Item(68, 11).Value2 = 960
$ws.Cells.Item(68, 12).Value2 = 2400
$ws.Cells.Item(68, 13).Value2 = -149
$ws.Cells.Item(68, 14).Value2 = -4022

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value2 = 533.3333
$ws.Cells.Item(71, 9).Value2 = 320
$ws.Cells.Item(71, 10).Value2 = 800
$ws.Cells.Item(71, 11).Value2 = 2880
$ws.Cells.Item(71, 12).Value2 = 7200
$ws.Cells.Item(71, 13).Value2 = 1176
$ws.Cells.Item(71, 14).Value2 = -15312

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value2 = 1389.9166
$ws.Cells.Item(92, 9).Value2 = 1362.6666
$ws.Cells.Item(92, 10).Value2 = 1417.1666
$ws.Cells.Item(92, 11).Value2 = 4087.9998
$ws.Cells.Item(92, 12).Value2 = 4251.4998
$ws.Cells.Item(92, 13).Value2 = -2839.9998
$ws.Cells.Item(92, 14).Value2 = -6747.4998

# GSM row 110
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(110, 8).Value2 = 48694
$ws.Cells.Item(110, 10).Value2 = 48694
$ws.Cells.Item(110, 12).Value2 = 48694
$ws.Cells.Item(110, 14).Value2 = -56874

# LTW row 81
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(81, 8).Value2 = 31090.5
$ws.Cells.Item(81, 10).Value2 = 31090.5
$ws.Cells.Item(81, 12).Value2 = 31090.5
$ws.Cells.Item(81, 14).Value2 = -33086.5

# LTW row 84
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(84, 8).Value2 = 31090.5
$ws.Cells.Item(84, 10).Value2 = 31090.5
$ws.Cells.Item(84, 12).Value2 = 93271.5
$ws.Cells.Item(84, 14).Value2 = -103255.5

# LTW row 111
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(111, 8).Value2 = 44387
$ws.Cells.Item(111, 10).Value2 = 44387
$ws.Cells.Item(111, 12).Value2 = 44387
$ws.Cells.Item(111, 14).Value2 = -52567

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value2 = 2739.8948
$ws.Cells.Item(136, 9).Value2 = 2273.2
$ws.Cells.Item(136, 10).Value2 = 4490
$ws.Cells.Item(136, 11).Value2 = 6819.599999999999
$ws.Cells.Item(136, 12).Value2 = 13470
$ws.Cells.Item(136, 13).Value2 = -4269.599999999999
$ws.Cells.Item(136, 14).Value2 = -18570

# WVR row 16
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value2 = 46523.5
$ws.Cells.Item(16, 10).Value2 = 46523.5
$ws.Cells.Item(16, 12).Value2 = 46523.5
$ws.Cells.Item(16, 14).Value2 = -47107.5

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value2 = 1335.6428
$ws.Cells.Item(81, 9).Value2 = 1299.909
$ws.Cells.Item(81, 10).Value2 = 1466.6666
$ws.Cells.Item(81, 11).Value2 = 2599.818
$ws.Cells.Item(81, 12).Value2 = 2933.3332
$ws.Cells.Item(81, 13).Value2 = -1538.818
$ws.Cells.Item(81, 14).Value2 = -5055.3332

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value2 = 1335.6428
$ws.Cells.Item(84, 9).Value2 = 1299.909
$ws.Cells.Item(84, 10).Value2 = 1466.6666
$ws.Cells.Item(84, 11).Value2 = 12999.09
$ws.Cells.Item(84, 12).Value2 = 14666.666
$ws.Cells.Item(84, 13).Value2 = -7695.09
$ws.Cells.Item(84, 14).Value2 = -25274.666

# WVR row 95
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value2 = 41172
$ws.Cells.Item(95, 10).Value2 = 41172
$ws.Cells.Item(95, 12).Value2 = 41172
$ws.Cells.Item(95, 14).Value2 = -46664

# WVR row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value2 = 47978
$ws.Cells.Item(119, 10).Value2 = 47978
$ws.Cells.Item(119, 12).Value2 = 47978
$ws.Cells.Item(119, 14).Value2 = -57654

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value2 = 40429
$ws.Cells.Item(123, 10).Value2 = 40429
$ws.Cells.Item(123, 12).Value2 = 40429
$ws.Cells.Item(123, 14).Value2 = -50229

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value2 = 16666.104
$ws.Cells.Item(136, 9).Value2 = 38170.184
$ws.Cells.Item(136, 10).Value2 = 2504.878
$ws.Cells.Item(136, 11).Value2 = 114510.552
$ws.Cells.Item(136, 12).Value2 = 7514.634
$ws.Cells.Item(136, 13).Value2 = -111960.552
$ws.Cells.Item(136, 14).Value2 = -12614.634
